$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 89
$ws.Range("I8").Value = 89
$ws.Range("K8").Value = 267
$ws.Range("M8").Value = -128

# Hunk 1: ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 403.10938
$ws.Range("I17").Value = 129.5
$ws.Range("J17").Value = 421.35
$ws.Range("K17").Value = 388.5
$ws.Range("L17").Value = 1264.05
$ws.Range("M17").Value = -220.5
$ws.Range("N17").Value = -1600.05

# Hunk 2: ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3445.3076
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 3865.4443
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 3865.4443
$ws.Range("M74").Value = -1564
$ws.Range("N74").Value = -5737.4443

# Hunk 3: ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3445.3076
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 3865.4443
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 19327.2215
$ws.Range("M77").Value = -7820
$ws.Range("N77").Value = -28687.2215

# Hunk 4: ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 792.2632
$ws.Range("I92").Value = 599.86664
$ws.Range("J92").Value = 1513.75
$ws.Range("K92").Value = 599.86664
$ws.Range("L92").Value = 1513.75
$ws.Range("M92").Value = 648.13336
$ws.Range("N92").Value = -4009.75

# Hunk 5: ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 9508.076999999999
$ws.Range("I98").Value = 9508.076999999999
$ws.Range("K98").Value = 9508.076999999999
$ws.Range("M98").Value = -8010.076999999999

# Hunk 6: ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1326.25
$ws.Range("I100").Value = 1390.5555
$ws.Range("J100").Value = 1133.3334
$ws.Range("K100").Value = 1390.5555
$ws.Range("L100").Value = 1133.3334
$ws.Range("M100").Value = -849.5554999999999
$ws.Range("N100").Value = -2215.3334

# Hunk 7: ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 552.85
$ws.Range("I107").Value = 621.63635
$ws.Range("J107").Value = 468.77777
$ws.Range("K107").Value = 621.63635
$ws.Range("L107").Value = 468.77777
$ws.Range("M107").Value = 1298.36365
$ws.Range("N107").Value = -4308.77777

# Hunk 8: ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1712.25
$ws.Range("I116").Value = 1499.5
$ws.Range("K116").Value = 1499.5
$ws.Range("M116").Value = 1942.5

# Hunk 9: ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 9508.076999999999
$ws.Range("I122").Value = 9508.076999999999
$ws.Range("K122").Value = 28524.231
$ws.Range("M122").Value = -26074.231

# Hunk 10: ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 999.8461
$ws.Range("I131").Value = 708.9091
$ws.Range("K131").Value = 2126.7273
$ws.Range("M131").Value = 2913.2727

# Hunk 11: ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3040588.8
$ws.Range("I132").Value = 3106623.2
$ws.Range("K132").Value = 9319869.600000001
$ws.Range("M132").Value = -9317339.600000001

# Hunk 12: ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1627.129
$ws.Range("I137").Value = 1373.381
$ws.Range("K137").Value = 4120.143
$ws.Range("M137").Value = -1570.143

# Hunk 13: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17121.305
$ws.Range("I32").Value = 18244.031
$ws.Range("J32").Value = 4771.3335
$ws.Range("K32").Value = 18244.031
$ws.Range("L32").Value = 4771.3335
$ws.Range("M32").Value = -17957.031
$ws.Range("N32").Value = -5345.3335

# Hunk 14: ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1709.931
$ws.Range("I74").Value = 1373.6296
$ws.Range("J74").Value = 6250
$ws.Range("K74").Value = 1373.6296
$ws.Range("L74").Value = 6250
$ws.Range("M74").Value = -499.6296
$ws.Range("N74").Value = -7998

# Hunk 15: ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1709.931
$ws.Range("I77").Value = 1373.6296
$ws.Range("J77").Value = 6250
$ws.Range("K77").Value = 6868.148
$ws.Range("L77").Value = 31250
$ws.Range("M77").Value = -2500.148
$ws.Range("N77").Value = -39986

# Hunk 16: ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4936.5
$ws.Range("I132").Value = 6190.3335
$ws.Range("K132").Value = 18571.0005
$ws.Range("M132").Value = -16041.0005

# Hunk 17: BSM row 10
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 19996
$ws.Range("J10").Value = 19996
$ws.Range("L10").Value = 19996
$ws.Range("N10").Value = -20276

# Hunk 18: BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3157.2104
$ws.Range("I20").Value = 3406.5833
$ws.Range("J20").Value = 2729.7144
$ws.Range("K20").Value = 3406.5833
$ws.Range("L20").Value = 2729.7144
$ws.Range("M20").Value = -3159.5833
$ws.Range("N20").Value = -3223.7144

# Hunk 19: BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1584.5555
$ws.Range("I86").Value = 1428.5714
$ws.Range("J86").Value = 1752.5385
$ws.Range("K86").Value = 1428.5714
$ws.Range("L86").Value = 1752.5385
$ws.Range("M86").Value = -305.5714
$ws.Range("N86").Value = -3998.5385

# Hunk 20: BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1584.5555
$ws.Range("I89").Value = 1428.5714
$ws.Range("J89").Value = 1752.5385
$ws.Range("K89").Value = 7142.857
$ws.Range("L89").Value = 8762.692500000001
$ws.Range("M89").Value = -1526.857
$ws.Range("N89").Value = -19994.6925

# Hunk 21: BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 625.36365
$ws.Range("I94").Value = 496.5
$ws.Range("J94").Value = 780
$ws.Range("K94").Value = 496.5
$ws.Range("L94").Value = 780
$ws.Range("M94").Value = -45.5
$ws.Range("N94").Value = -1682

# Hunk 22: CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5407693
$ws.Range("I31").Value = 2289.743
$ws.Range("J31").Value = 100002250
$ws.Range("K31").Value = 2289.743
$ws.Range("L31").Value = 100002250
$ws.Range("M31").Value = -1994.743
$ws.Range("N31").Value = -100002840

# Hunk 23: CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5407693
$ws.Range("I34").Value = 2289.743
$ws.Range("J34").Value = 100002250
$ws.Range("K34").Value = 2289.743
$ws.Range("L34").Value = 100002250
$ws.Range("M34").Value = -2087.743
$ws.Range("N34").Value = -100002654

# Hunk 24: CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 918.6667
$ws.Range("I122").Value = 878
$ws.Range("K122").Value = 2634
$ws.Range("M122").Value = -184

# Hunk 25: CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1285.8334
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""

# Hunk 26: CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 556950.4399999999
$ws.Range("J122").Value = 1235823.1
$ws.Range("L122").Value = 11122407.9
$ws.Range("N122").Value = -11127307.9

# Hunk 27: CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1285.8334
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").Value = ""

# Hunk 28: GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1983.3334
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 1983.3334
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 1983.3334
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = -5227.3334

# Hunk 29: GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 41668250
$ws.Range("J113").Value = 2166.3333
$ws.Range("L113").Value = 2166.3333
$ws.Range("N113").Value = -6506.3333

# Hunk 30: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 64936.22
$ws.Range("I132").Value = 76109.664
$ws.Range("J132").Value = 4599.6
$ws.Range("K132").Value = 228328.992
$ws.Range("L132").Value = 13798.8
$ws.Range("M132").Value = -225798.992
$ws.Range("N132").Value = -18858.8

# Hunk 31: WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 457.7647
$ws.Range("I113").Value = 562.8570999999999
$ws.Range("J113").Value = 384.2
$ws.Range("K113").Value = 1688.5713
$ws.Range("L113").Value = 1152.6
$ws.Range("M113").Value = 481.4287000000002
$ws.Range("N113").Value = -5492.6
